$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 2 'Bitcoin'
Set-TextCell 2 3 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
Set-TextCell 2 4 '29.367.23'
Set-TextCell 2 5 '  -0.40%  '
Set-TextCell 3 2 'Ethereum'
Set-TextCell 3 3 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
Set-TextCell 3 4 '1.845.82'
Set-TextCell 3 5 '  -0.22%  '
Set-TextCell 4 2 'TetherUSD'
Set-TextCell 4 3 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
Set-TextCell 4 4 '0.9987'
Set-TextCell 4 5 '  +0.03%  '
Set-TextCell 5 2 'BNB'
Set-TextCell 5 3 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextCell 5 4 '240.09'
Set-TextCell 5 5 '  -0.79%  '
Set-TextCell 6 2 'XRP'
Set-TextCell 6 3 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextCell 6 4 '0.6304'
Set-TextCell 6 5 '  +0.39%  '
Set-TextCell 7 2 'USDC'
Set-TextCell 7 3 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextCell 7 4 '1.000'
Set-TextCell 7 5 '  +0.05%  '
Set-TextCell 8 2 'Dogecoin'
Set-TextCell 8 3 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell 8 4 '0.07538'
Set-TextCell 8 5 '  +0.05%  '
Set-TextCell 9 2 'Cardano'
Set-TextCell 9 3 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell 9 4 '0.2956'
Set-TextCell 9 5 '  -0.72%  '
Set-TextCell 10 2 'Solana'
Set-TextCell 10 3 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell 10 4 '24.43'
Set-TextCell 10 5 '  +0.51%  '
Set-TextCell 11 2 'TRON'
Set-TextCell 11 3 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextCell 11 4 '0.07714'
Set-TextCell 11 5 '  +0.15%  '
Set-TextCell 12 2 'WrappedEther'
Set-TextCell 12 3 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 12 4 '1.842.45'
Set-TextCell 12 5 '  -4.10%  '
Set-TextCell 13 2 'Polkadot'
Set-TextCell 13 3 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 13 4 '4.992'
Set-TextCell 13 5 '  -0.22%  '
Set-TextCell 14 2 'Polygon'
Set-TextCell 14 3 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell 14 4 '0.6831'
Set-TextCell 14 5 '  -0.35%  '
Set-TextCell 15 2 'ShibaInu'
Set-TextCell 15 3 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell 15 4 '0.00001000'
Set-TextCell 15 5 '  +2.25%  '
Set-TextCell 16 2 'Litecoin'
Set-TextCell 16 3 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 16 4 '82.82'
Set-TextCell 16 5 '  -1.07%  '
Set-TextCell 17 2 'Uniswap'
Set-TextCell 17 3 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell 17 4 '6.128'
Set-TextCell 17 5 '  -1.59%  '
Set-TextCell 18 2 'WrappedBTC'
Set-TextCell 18 3 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextCell 18 4 '29.398.33'
Set-TextCell 18 5 '  -0.63%  '
Set-TextCell 19 2 'BitcoinCash'
Set-TextCell 19 3 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextCell 19 4 '227.74'
Set-TextCell 19 5 '  -2.73%  '
Set-TextCell 20 2 'Avalanche'
Set-TextCell 20 3 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell 20 4 '12.42'
Set-TextCell 20 5 '  -0.59%  '
Set-TextCell 21 2 'Dai'
Set-TextCell 21 3 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell 21 4 '0.9998'
Set-TextCell 21 5 '  -0.11%  '
Set-TextCell 22 2 'Chainlink'
Set-TextCell 22 3 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextCell 22 4 '7.534'
Set-TextCell 22 5 '  -1.19%  '
Set-TextCell 23 2 'LEO'
Set-TextCell 23 3 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell 23 4 '3.962'
Set-TextCell 23 5 '  -0.47%  '
Set-TextCell 24 2 'BinanceUSD'
Set-TextCell 24 3 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell 24 4 '1.001'
Set-TextCell 24 5 '  +0.14%  '
Set-TextCell 25 2 'RocketPoolETH'
Set-TextCell 25 3 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextCell 25 4 '5.956.16'
Set-TextCell 25 5 '  +179.85%  '
Set-TextCell 26 2 'Monero'
Set-TextCell 26 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 26 4 '157.37'
Set-TextCell 26 5 '  +1.51%  '
Set-TextCell 27 2 'Stellar'
Set-TextCell 27 3 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 27 4 '0.1396'
Set-TextCell 27 5 '  +0.19%  '
Set-TextCell 28 2 'Cosmos'
Set-TextCell 28 3 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 28 4 '8.368'
Set-TextCell 28 5 '  -0.80%  '
Set-TextCell 29 2 'EthereumClassic'
Set-TextCell 29 3 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell 29 4 '17.64'
Set-TextCell 29 5 '  -0.50%  '
Set-TextCell 30 2 'PancakeSwap'
Set-TextCell 30 3 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextCell 30 4 '1.464'
Set-TextCell 30 5 '  -0.93%  '
Set-TextCell 31 2 'Toncoin'
Set-TextCell 31 3 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 31 4 '1.255'
Set-TextCell 31 5 '  -0.60%  '
Set-TextCell 32 2 'Hedera'
Set-TextCell 32 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 32 4 '0.05676'
Set-TextCell 32 5 '  -3.04%  '
Set-TextCell 33 2 'Filecoin'
Set-TextCell 33 3 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 33 4 '4.124'
Set-TextCell 33 5 '  +0.64%  '
Set-TextCell 34 2 'InternetComputer(DFINITY)'
Set-TextCell 34 3 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 34 4 '4.020'
Set-TextCell 34 5 '  -0.52%  '
Set-TextCell 35 2 'LidoDAOToken'
Set-TextCell 35 3 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell 35 4 '1.843'
Set-TextCell 35 5 '  -3.14%  '
Set-TextCell 36 2 'ARBITRUM'
Set-TextCell 36 3 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 36 4 '1.154'
Set-TextCell 36 5 '  -1.32%  '
Set-TextCell 37 2 'ImmutableX'
Set-TextCell 37 3 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 37 4 '0.7132'
Set-TextCell 37 5 '  -1.42%  '
Set-TextCell 38 2 'HuobiToken'
Set-TextCell 38 3 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 38 4 '2.591'
Set-TextCell 38 5 '  +0.16%  '
Set-TextCell 39 2 'Maker'
Set-TextCell 39 3 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell 39 4 '1.258.61'
Set-TextCell 39 5 '  +1.50%  '
Set-TextCell 40 2 'VeChain'
Set-TextCell 40 3 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 40 4 '0.01813'
Set-TextCell 40 5 '  +1.31%  '
Set-TextCell 41 2 'MXToken'
Set-TextCell 41 3 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell 41 4 '2.782'
Set-TextCell 41 5 '  -0.40%  '
Set-TextCell 42 2 'TrustWalletToken'
Set-TextCell 42 3 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 42 4 '0.9121'
Set-TextCell 42 5 '  +0.45%  '
Set-TextCell 43 2 'FraxShare'
Set-TextCell 43 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 43 4 '6.215'
Set-TextCell 43 5 '  +1.16%  '
Set-TextCell 44 2 'PaxDollar'
Set-TextCell 44 3 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextCell 44 4 '1.000'
Set-TextCell 44 5 '  +0.06%  '
Set-TextCell 45 2 'Quant'
Set-TextCell 45 3 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextCell 45 4 '101.02'
Set-TextCell 45 5 '  -1.09%  '
Set-TextCell 46 2 'Aave'
Set-TextCell 46 3 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell 46 4 '66.30'
Set-TextCell 46 5 '  -1.17%  '
Set-TextCell 47 2 'BabyDogeCoin'
Set-TextCell 47 3 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell 47 4 '0.00000000118'
Set-TextCell 47 5 '  +0.29%  '
Set-TextCell 48 2 'Aptos'
Set-TextCell 48 3 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 48 4 '7.036'
Set-TextCell 48 5 '  -3.97%  '
Set-TextCell 49 2 'TheSandbox'
Set-TextCell 49 3 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextCell 49 4 '0.4021'
Set-TextCell 49 5 '  -0.32%  '
Set-TextCell 50 2 'EnergySwap'
Set-TextCell 50 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 50 4 '9.083'
Set-TextCell 50 5 '  -0.64%  '
Set-TextCell 51 2 'RenderToken'
Set-TextCell 51 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 51 4 '1.687'
Set-TextCell 51 5 '  -2.34%  '
